$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Zeitachse"
$ws1.PageSetup.PrintArea = "`$A`$1:`$N`$17"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C4").Value = "Mockups erstellen"
$ws2.Range("E4").ClearContents()
